$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.357.26'
$ws.Range("E2").Value = '  +0.55%  '
$ws.Range("D3").Value = '1.597.38'
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '211.98'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").Value = '0.500'
$ws.Range("E6").Value = '  -0.28%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").Value = '19.15'
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("D12").Value = '1.821.47'
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("D13").Value = '1.606.97'
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("E15").Value = '  -1.18%  '
$ws.Range("D16").Value = '63.47'
$ws.Range("E16").Value = '  -0.38%  '
$ws.Range("D17").Value = '26.328.68'
$ws.Range("E17").Value = '  +0.51%  '
$ws.Range("D18").Value = '229.61'
$ws.Range("E18").Value = '  +7.10%  '
$ws.Range("D19").Value = '7.67'
$ws.Range("E19").Value = '  +4.22%  '
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("D21").Value = '0.999'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("D22").Value = '4.25'
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("E23").Value = '  +2.80%  '
$ws.Range("E24").Value = '  -1.31%  '
$ws.Range("D25").Value = '146.57'
$ws.Range("E25").Value = '  +1.13%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").Value = '6.98'
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("E28").Value = '  +0.56%  '
$ws.Range("D29").Value = '15.37'
$ws.Range("E29").Value = '  +1.81%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("D32").Value = '1.498.66'
$ws.Range("E32").Value = '  +4.86%  '
$ws.Range("E33").Value = '  +1.06%  '
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("E36").Value = '  +0.37%  '
$ws.Range("D37").Value = '0.571'
$ws.Range("E37").Value = '  -2.85%  '
$ws.Range("E38").Value = '  -0.84%  '
$ws.Range("E39").Value = '  -0.85%  '
$ws.Range("E40").Value = '  -2.06%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").Value = '0.940'
$ws.Range("E42").Value = '  -4.28%  '
$ws.Range("E43").Value = '  +1.69%  '
$ws.Range("D44").Value = '1.733.88'
$ws.Range("E44").Value = '  +0.37%  '
$ws.Range("D45").Value = '0.758'
$ws.Range("E45").Value = '  -1.00%  '
$ws.Range("D46").Value = '60.68'
$ws.Range("E46").Value = '  -0.66%  '
$ws.Range("D47").Value = '88.41'
$ws.Range("E47").Value = '  +1.57%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0103'
$ws.Range("E48").Value = '  +0.16%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '1.48'
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.0501'
$ws.Range("E50").Value = '  -0.17%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '0.0959'
$ws.Range("E51").Value = '  +0.20%  '
